$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "./model_output/2025-08-12-19-13-56-None"
$ws.Range("B6").Value = 0.5292857142857142
